$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.371.29'
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').Value = '1.907.82'
$ws.Range('E3').Value = '  +2.61%  '
$ws.Range('E4').Value = '  -0.57%  '
$ws.Range('D5').Value = '''246.29'
$ws.Range('E5').Value = '  +2.92%  '
$ws.Range('D6').Value = '''0.666'
$ws.Range('E6').Value = '  +6.95%  '
$ws.Range('E7').Value = '  -0.53%  '
$ws.Range('D8').Value = '''41.48'
$ws.Range('E8').Value = '  -1.97%  '
$ws.Range('E9').Value = '  +5.93%  '
$ws.Range('D10').Value = '''52.91'
$ws.Range('E10').Value = '  +12.74%  '
$ws.Range('E11').Value = '  +4.06%  '
$ws.Range('E12').Value = '  +0.49%  '
$ws.Range('E13').Value = '  +2.36%  '
$ws.Range('D14').Value = '''12.09'
$ws.Range('E14').Value = '  +5.29%  '
$ws.Range('D15').Value = '''0.700'
$ws.Range('E15').Value = '  +3.55%  '
$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').Value = '''4.88'
$ws.Range('E16').Value = '  +3.63%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '1.899.53'
$ws.Range('E17').Value = '  +2.61%  '
$ws.Range('D18').Value = '35.357.85'
$ws.Range('E18').Value = '  +0.57%  '
$ws.Range('D19').Value = '''72.30'
$ws.Range('E19').Value = '  +3.43%  '
$ws.Range('D20').Value = '0.0₃0830'
$ws.Range('E20').Value = '  +4.53%  '
$ws.Range('D21').Value = '''240.13'
$ws.Range('E21').Value = '  -0.43%  '
$ws.Range('D22').Value = '''12.54'
$ws.Range('E22').Value = '  +2.54%  '
$ws.Range('E23').Value = '  +2.23%  '
$ws.Range('E24').Value = '  -0.50%  '
$ws.Range('D25').Value = '''2.29'
$ws.Range('E25').Value = '  +1.01%  '
$ws.Range('D26').Value = '''2.33'
$ws.Range('E26').Value = '  +23.75%  '
$ws.Range('D28').Value = '''8.44'
$ws.Range('E28').Value = '  +4.89%  '
$ws.Range('D29').Value = '''18.46'
$ws.Range('E29').Value = '  +4.35%  '
$ws.Range('D30').Value = '''0.128'
$ws.Range('E30').Value = '  +2.87%  '
$ws.Range('E31').Value = '  +3.48%  '
$ws.Range('D32').Value = '''0.0567'
$ws.Range('E32').Value = '  +0.97%  '
$ws.Range('D33').Value = '''1.01'
$ws.Range('E33').Value = '  -0.44%  '
$ws.Range('D34').Value = '''0.931'
$ws.Range('E34').Value = '  +14.66%  '
$ws.Range('E35').Value = '  +2.11%  '
$ws.Range('E36').Value = '  -4.24%  '
$ws.Range('E37').Value = '  -0.12%  '
$ws.Range('E38').Value = '  +2.15%  '
$ws.Range('E39').Value = '  +0.78%  '
$ws.Range('E40').Value = '  +3.71%  '
$ws.Range('D41').Value = '''16.35'
$ws.Range('E41').Value = '  +8.70%  '
$ws.Range('D42').Value = '''0.0635'
$ws.Range('E42').Value = '  +7.40%  '
$ws.Range('D43').Value = '''90.12'
$ws.Range('E43').Value = '  +0.12%  '
$ws.Range('D44').Value = '1.341.56'
$ws.Range('D45').Value = '''2.40'
$ws.Range('E45').Value = '  +3.17%  '
$ws.Range('D46').Value = '''48.10'
$ws.Range('E46').Value = '  +39.07%  '
$ws.Range('E47').Value = '  +1.58%  '
$ws.Range('E48').Value = '  -0.98%  '
$ws.Range('E49').Value = '  -0.12%  '
$ws.Range('B50').Value = 'Gas'
$ws.Range('C50').Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range('D50').Value = '''11.64'
$ws.Range('E50').Value = '  -4.80%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.091.96'
$ws.Range('E51').Value = '  +2.41%  '
